$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay as text (e.g. trailing
# zeros like "0.610" or multi-dot thousands separators like "45.431.70")
# -- force text format before assigning so Excel does not coerce them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.431.70'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.371.59'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.03'
$ws.Range("E5").Value = '  -1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.21'
$ws.Range("E6").Value = '  -3.49%  '
$ws.Range("E7").Value = '  -1.34%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.610'
$ws.Range("E9").Value = '  -3.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.79'
$ws.Range("E10").Value = '  -3.69%  '
$ws.Range("E11").Value = '  -1.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.49'
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.981'
$ws.Range("E14").Value = '  -3.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.733.11'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.33'
$ws.Range("E16").Value = '  -3.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.381.12'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.465.28'
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.89'
$ws.Range("E19").Value = '  +5.50%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000107'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.25'
$ws.Range("E21").Value = '  -5.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.34'
$ws.Range("E22").Value = '  -2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.54'
$ws.Range("E23").Value = '  -0.86%  '
$ws.Range("E24").Value = '  -3.80%  '
$ws.Range("E25").Value = '  +2.27%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.15'
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.20'
$ws.Range("E28").Value = '  -6.31%  '
$ws.Range("E29").Value = '  -1.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0980'
$ws.Range("E30").Value = '  +4.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.35'
$ws.Range("E31").Value = '  -2.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.13'
$ws.Range("E32").Value = '  -6.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '166.73'
$ws.Range("E33").Value = '  -1.57%  '
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.130'
$ws.Range("E35").Value = '  -2.11%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.118'
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.70'
$ws.Range("E37").Value = '  -2.92%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.01'
$ws.Range("E38").Value = '  +1.81%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.89'
$ws.Range("E39").Value = '  +8.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.95'
$ws.Range("E40").Value = '  -3.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0355'
$ws.Range("E41").Value = '  -3.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.57'
$ws.Range("E42").Value = '  -6.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.96'
$ws.Range("E43").Value = '  -2.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.228'
$ws.Range("E44").Value = '  -5.49%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.67'
$ws.Range("E46").Value = '  -8.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.823.55'
$ws.Range("E47").Value = '  +9.75%  '
$ws.Range("E48").Value = '  +5.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.79'
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.28'
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.90'
$ws.Range("E51").Value = '  -6.88%  '
